$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (dates 2021-08-24 .. 2021-09-01) to append below the
# existing data (which previously ended at row 357 / 2021-08-23).
$data = @(
    @(44432, 0, 4, 101.7293997965412),
    @(44433, 0, 4, 101.7293997965412),
    @(44434, 0, 4, 101.7293997965412),
    @(44435, 0, 2, 50.8646998982706),
    @(44436, 0, 2, 50.8646998982706),
    @(44437, 0, 0, 0),
    @(44438, 0, 0, 0),
    @(44439, 0, 0, 0),
    @(44440, 0, 0, 0)
)

$startRow = 358
$lastRow = $startRow + $data.Length - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Column A carries the date style (bold, thin box border, centered/top
# alignment, custom date number format) used throughout the sheet. Copy
# that formatting from the last pre-existing data row (357) down onto
# the new date cells so the appended rows match the rest of the column.
$ws.Range("A357").Copy() | Out-Null
$ws.Range("A$startRow`:A$lastRow").PasteSpecial(-4122) | Out-Null
